# Auto-generated Excel COM-interop script to apply numeric value updates
# to the Zodiark_Profits workbook (8 job-sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 1245.591
$ws.Range("J17").Value = 1270.1111
$ws.Range("L17").Value = 3810.3333
$ws.Range("N17").Value = -4146.3333
$ws.Range("H70").Value = 1180.7037
$ws.Range("I70").Value = 1212.7
$ws.Range("J70").Value = 1089.2858
$ws.Range("K70").Value = 3638.1
$ws.Range("L70").Value = 3267.8574
$ws.Range("M70").Value = -3368.1
$ws.Range("N70").Value = -3807.8574
$ws.Range("H73").Value = 1180.7037
$ws.Range("I73").Value = 1212.7
$ws.Range("J73").Value = 1089.2858
$ws.Range("K73").Value = 3638.1
$ws.Range("L73").Value = 3267.8574
$ws.Range("M73").Value = -2702.1
$ws.Range("N73").Value = -5139.857400000001
$ws.Range("H76").Value = 4900.5454
$ws.Range("J76").Value = 5981.4
$ws.Range("L76").Value = 5981.4
$ws.Range("N76").Value = -6611.4
$ws.Range("H79").Value = 4900.5454
$ws.Range("J79").Value = 5981.4
$ws.Range("L79").Value = 5981.4
$ws.Range("N79").Value = -8165.4
$ws.Range("H107").Value = 2090
$ws.Range("J107").Value = 2275.7144
$ws.Range("L107").Value = 2275.7144
$ws.Range("N107").Value = -6115.7144
$ws.Range("H132").Value = 2674.6
$ws.Range("I132").Value = 2574.7307
$ws.Range("J132").Value = 3323.75
$ws.Range("K132").Value = 7724.1921
$ws.Range("L132").Value = 9971.25
$ws.Range("M132").Value = -5194.1921
$ws.Range("N132").Value = -15031.25
$ws.Range("H137").Value = 2155.3333
$ws.Range("I137").Value = 2470
$ws.Range("K137").Value = 7410
$ws.Range("M137").Value = -4860
$ws.Range("H138").Value = 10206162
$ws.Range("J138").Value = 2516.1045
$ws.Range("L138").Value = 7548.3135
$ws.Range("N138").Value = -17828.3135
$ws.Range("H141").Value = 3691.375
$ws.Range("I141").Value = 2647.2856
$ws.Range("K141").Value = 7941.8568
$ws.Range("M141").Value = -2761.8568

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 2516.3381
$ws.Range("I32").Value = 1562.8334
$ws.Range("K32").Value = 1562.8334
$ws.Range("M32").Value = -1275.8334
$ws.Range("H45").Value = 1741.75
$ws.Range("I45").Value = 2343
$ws.Range("K45").Value = 2343
$ws.Range("M45").Value = -1966
$ws.Range("H48").Value = 118997.5
$ws.Range("J48").Value = 118997.5
$ws.Range("L48").Value = 118997.5
$ws.Range("N48").Value = -119765.5
$ws.Range("H61").Value = 1629
$ws.Range("I61").Value = 1321.6957
$ws.Range("K61").Value = 1321.6957
$ws.Range("M61").Value = -1109.6957
$ws.Range("H136").Value = 1629
$ws.Range("I136").Value = 1321.6957
$ws.Range("K136").Value = 3965.0871
$ws.Range("M136").Value = -1415.0871

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H41").Value = 119965
$ws.Range("J41").Value = 119965
$ws.Range("L41").Value = 119965
$ws.Range("N41").Value = -120741
$ws.Range("H42").Value = 119952.5
$ws.Range("J42").Value = 119952.5
$ws.Range("L42").Value = 119952.5
$ws.Range("N42").Value = -120608.5
$ws.Range("H43").Value = 171641.67
$ws.Range("J43").Value = 171641.67
$ws.Range("L43").Value = 171641.67
$ws.Range("N43").Value = -172003.67
$ws.Range("H47").Value = 746648.3
$ws.Range("J47").Value = 746648.3
$ws.Range("L47").Value = 746648.3
$ws.Range("N47").Value = -747688.3
$ws.Range("H48").Value = 119975
$ws.Range("J48").Value = 119975
$ws.Range("L48").Value = 119975
$ws.Range("N48").Value = -120805

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H62").Value = 12284.904
$ws.Range("I62").Value = 9124.5
$ws.Range("K62").Value = 9124.5
$ws.Range("M62").Value = -8500.5
$ws.Range("H65").Value = 12284.904
$ws.Range("I65").Value = 9124.5
$ws.Range("K65").Value = 45622.5
$ws.Range("M65").Value = -42502.5
$ws.Range("H86").Value = 76928250
$ws.Range("I86").Value = 200003200
$ws.Range("K86").Value = 200003200
$ws.Range("M86").Value = -200002077
$ws.Range("H89").Value = 76928250
$ws.Range("I89").Value = 200003200
$ws.Range("K89").Value = 1000016000
$ws.Range("M89").Value = -1000010384
$ws.Range("H105").Value = 20612.223
$ws.Range("I105").Value = 26104.285
$ws.Range("K105").Value = 26104.285
$ws.Range("M105").Value = -24357.285
$ws.Range("H122").Value = 2635
$ws.Range("I122").Value = 2307.4443
$ws.Range("K122").Value = 6922.3329
$ws.Range("M122").Value = -4472.3329

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 15151738
$ws.Range("I2").Value = 187.22223
$ws.Range("J2").Value = 25641274
$ws.Range("K2").Value = 1123.33338
$ws.Range("L2").Value = 153847644
$ws.Range("M2").Value = -1010.33338
$ws.Range("N2").Value = -153847870
$ws.Range("H12").Value = 102.42857
$ws.Range("I12").Value = 173.5
$ws.Range("J12").Value = 74
$ws.Range("K12").Value = 520.5
$ws.Range("L12").Value = 222
$ws.Range("M12").Value = -347.5
$ws.Range("N12").Value = -568
$ws.Range("H121").Value = 3287.7778
$ws.Range("J121").Value = 3451.7646
$ws.Range("L121").Value = 10355.2938
$ws.Range("N121").Value = -12975.2938
$ws.Range("H137").Value = 2392.6667
$ws.Range("I137").Value = 2472.5715
$ws.Range("K137").Value = 7417.7145
$ws.Range("M137").Value = -2317.7145

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 3187.3333
$ws.Range("I102").Value = 3187.3333
$ws.Range("K102").Value = 3187.3333
$ws.Range("M102").Value = -1565.3333
$ws.Range("H126").Value = 4097.909
$ws.Range("I126").Value = 3507.7
$ws.Range("K126").Value = 10523.1
$ws.Range("M126").Value = -8053.099999999999
$ws.Range("H132").Value = 2361.238
$ws.Range("I132").Value = 2361.238
$ws.Range("K132").Value = 7083.714
$ws.Range("M132").Value = -4553.714

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 6079.5
$ws.Range("I7").Value = 3466.3333
$ws.Range("J7").Value = 9999.25
$ws.Range("K7").Value = 3466.3333
$ws.Range("L7").Value = 9999.25
$ws.Range("M7").Value = -3354.3333
$ws.Range("N7").Value = -10223.25
$ws.Range("H16").Value = 1036.0834
$ws.Range("I16").Value = 1298.4445
$ws.Range("K16").Value = 1298.4445
$ws.Range("M16").Value = -1128.4445
$ws.Range("H61").Value = 1022.13336
$ws.Range("I61").Value = 931.6923
$ws.Range("K61").Value = 931.6923
$ws.Range("M61").Value = -729.6923
$ws.Range("H82").Value = 1619.8649
$ws.Range("I82").Value = 908
$ws.Range("K82").Value = 908
$ws.Range("M82").Value = -547
$ws.Range("H85").Value = 1619.8649
$ws.Range("I85").Value = 908
$ws.Range("K85").Value = 908
$ws.Range("M85").Value = 340
$ws.Range("H113").Value = 1022.13336
$ws.Range("I113").Value = 931.6923
$ws.Range("K113").Value = 931.6923
$ws.Range("M113").Value = 1238.3077
$ws.Range("H126").Value = 6079.5
$ws.Range("I126").Value = 3466.3333
$ws.Range("J126").Value = 9999.25
$ws.Range("K126").Value = 10398.9999
$ws.Range("L126").Value = 29997.75
$ws.Range("M126").Value = -7928.999899999999
$ws.Range("N126").Value = -34937.75

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H41").Value = 19884.5
$ws.Range("I41").Value = 21493
$ws.Range("J41").Value = 18276
$ws.Range("K41").Value = 21493
$ws.Range("L41").Value = 18276
$ws.Range("M41").Value = -21103
$ws.Range("N41").Value = -19056
$ws.Range("H81").Value = 3961
$ws.Range("I81").Value = 2021.4286
$ws.Range("J81").Value = 10749.5
$ws.Range("K81").Value = 4042.8572
$ws.Range("L81").Value = 21499
$ws.Range("M81").Value = -2981.8572
$ws.Range("N81").Value = -23621
$ws.Range("H84").Value = 3961
$ws.Range("I84").Value = 2021.4286
$ws.Range("J84").Value = 10749.5
$ws.Range("K84").Value = 20214.286
$ws.Range("L84").Value = 107495
$ws.Range("M84").Value = -14910.286
$ws.Range("N84").Value = -118103
